# Update "想去人数" (number of people wanting to attend) counts for two
# events that appear both on the "展览" sheet and the "全部类型" sheet.
#   南宁·熊谷M动漫嘉年华（免费）      1332 -> 1334
#   南宁·第二届北极光动漫展          2849 -> 2852

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1334
$wsExhibit.Range("F3").Value = 2852

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1334
$wsAll.Range("F4").Value = 2852
